# Append a new "Filtered Feeds" row describing the latest 360Dx "Top Five
# Articles" roundup (picked up by the scraping workflow).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has data in rows 2..70 (row 1 is the header), so the
# new entry goes into row 71.
$newRow = 71

$newLink    = "https://www.360dx.com/business-news/top-five-articles-360dx-last-week-roche-vaginitis-test-fda-reclassification-cdx"
$newKeyword = "CDx"
$newTitle   = "Top Five Articles on 360Dx Last Week: Roche Vaginitis Test; FDA Reclassification of CDx Assays; More"

$ws.Cells.Item($newRow, 1).Value = $newLink
$ws.Cells.Item($newRow, 2).Value = $newKeyword
$ws.Cells.Item($newRow, 3).Value = $newTitle

# Turn the link cell into a real hyperlink, matching the styling already
# used by the other "link" cells in column A.
$ws.Hyperlinks.Add($ws.Cells.Item($newRow, 1), $newLink)
$ws.Cells.Item($newRow, 1).Style = $ws.Cells.Item($newRow - 1, 1).Style
